$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 253-255 with revised figures ---
$ws.Range("B253").Value2 = 5817095330000
$ws.Range("D253").Value2 = 257940867016.524

$ws.Range("B254").Value2 = 5950864520000
$ws.Range("D254").Value2 = 271913974347.8439

$ws.Range("B255").Value2 = 6045092150000
$ws.Range("D255").Value2 = 272103499031.1103

# --- Append new rows 256-258, copying the date-column style from row 255 ---
$ws.Range("A255").Copy($ws.Range("A256"))
$ws.Range("A256").Value2 = 44986
$ws.Range("B256").Value2 = 6077620130000
$ws.Range("C256").Value2 = 0.04620292787953972
$ws.Range("D256").Value2 = 280803844545.6288

$ws.Range("A255").Copy($ws.Range("A257"))
$ws.Range("A257").Value2 = 45017
$ws.Range("B257").Value2 = 6141246740000
$ws.Range("C257").Value2 = 0.04682075376731491
$ws.Range("D257").Value2 = 287537801437.8654

$ws.Range("A255").Copy($ws.Range("A258"))
$ws.Range("A258").Value2 = 45047
$ws.Range("B258").Value2 = 6224272840000
$ws.Range("C258").Value2 = 0.04506372010022171
$ws.Range("D258").Value2 = 280488889089.1721
